# Updated cryptos list on Fri Nov 10 21:41:25 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    # Force the cell to be treated as plain text so numeric-looking
    # strings (e.g. "252.24") are not coerced into floating point numbers,
    # then clear the temporary number-format style so the cell keeps its
    # original (default) style like the rest of the sheet.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "37.397.70"
Set-TextValue "E2" "  +2.09%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.099.32"
Set-TextValue "E3" "  +1.00%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.16%  "

# Row 5 - BNB
Set-TextValue "D5" "252.24"
Set-TextValue "E5" "  +1.70%  "

# Row 6 - XRP
Set-TextValue "D6" "0.668"
Set-TextValue "E6" "  +0.25%  "

# Row 7 - USDC
Set-TextValue "E7" "  -0.06%  "

# Row 8 - Solana
Set-TextValue "D8" "54.27"
Set-TextValue "E8" "  +19.90%  "

# Row 9 - OKB
Set-TextValue "D9" "62.48"
Set-TextValue "E9" "  +2.85%  "

# Row 10 - Cardano
Set-TextValue "E10" "  +4.21%  "

# Row 11 - Dogecoin
Set-TextValue "E11" "  +4.86%  "

# Row 12 - TRON
Set-TextValue "E12" "  +7.68%  "

# Row 13 - Chainlink
Set-TextValue "D13" "15.13"
Set-TextValue "E13" "  +3.81%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "2.400.20"
Set-TextValue "E14" "  +1.34%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.837"
Set-TextValue "E15" "  +2.52%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.096.60"
Set-TextValue "E16" "  +0.99%  "

# Row 17 - Polkadot
Set-TextValue "E17" "  +6.67%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "37.295.33"
Set-TextValue "E18" "  +1.86%  "

# Row 19 - Litecoin
Set-TextValue "D19" "73.27"
Set-TextValue "E19" "  +2.38%  "

# Row 20 - Avalanche
Set-TextValue "D20" "14.51"
Set-TextValue "E20" "  +14.62%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0855"
Set-TextValue "E21" "  +5.12%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "241.62"
Set-TextValue "E22" "  +1.53%  "

# Row 23 - Uniswap
Set-TextValue "E23" "  +6.30%  "

# Row 24 - Dai
Set-TextValue "E24" "  +0.00%  "

# Row 25 - Toncoin
Set-TextValue "E25" "  +0.82%  "

# Row 26 - Monero
Set-TextValue "D26" "171.83"
Set-TextValue "E26" "  +1.40%  "

# Row 27 - Cosmos
Set-TextValue "D27" "9.28"
Set-TextValue "E27" "  +5.43%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "20.96"
Set-TextValue "E28" "  +3.58%  "

# Row 29 - PancakeSwap
Set-TextValue "E29" "  +4.61%  "

# Row 30 - Stellar
Set-TextValue "E30" "  +2.12%  "

# Row 31 - was ImmutableX, now Gas (rows 31/32 swapped)
Set-TextValue "B31" "Gas"
Set-TextValue "C31" "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
Set-TextValue "D31" "22.77"
Set-TextValue "E31" "  +4.92%  "

# Row 32 - was Gas, now ImmutableX
Set-TextValue "B32" "ImmutableX"
Set-TextValue "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "1.10"
Set-TextValue "E32" "  +25.17%  "

# Row 33 - Filecoin
Set-TextValue "D33" "4.55"
Set-TextValue "E33" "  +4.19%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0625"
Set-TextValue "E34" "  +7.36%  "

# Row 35 - Kaspa
Set-TextValue "D35" "0.0906"
Set-TextValue "E35" "  +0.51%  "

# Row 36 - InternetComputer(DFINITY)
Set-TextValue "D36" "4.24"
Set-TextValue "E36" "  +6.40%  "

# Row 37 - BinanceUSD
Set-TextValue "E37" "  -0.05%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "2.27"
Set-TextValue "E38" "  -0.12%  "

# Row 39 - WEMIXToken
Set-TextValue "E39" "  -2.73%  "

# Row 40 - FTXToken
Set-TextValue "D40" "5.08"
Set-TextValue "E40" "  +148.67%  "

# Row 41 - TrustWalletToken
Set-TextValue "E41" "  +2.56%  "

# Row 42 - InjectiveProtocol
Set-TextValue "D42" "18.19"
Set-TextValue "E42" "  +12.70%  "

# Row 43 - VeChain
Set-TextValue "E43" "  +6.64%  "

# Row 44 - ARBITRUM
Set-TextValue "E44" "  +3.33%  "

# Row 45 - Cronos
Set-TextValue "D45" "0.0973"
Set-TextValue "E45" "  +19.00%  "

# Row 46 - Aave
Set-TextValue "D46" "99.50"
Set-TextValue "E46" "  +2.76%  "

# Row 47 - HuobiToken
Set-TextValue "D47" "2.81"
Set-TextValue "E47" "  +0.88%  "

# Row 48 - Maker
Set-TextValue "D48" "1.334.11"
Set-TextValue "E48" "  +1.05%  "

# Row 49 - MXToken
Set-TextValue "E49" "  +4.34%  "

# Row 50 - RenderToken
Set-TextValue "D50" "2.37"
Set-TextValue "E50" "  +7.45%  "

# Row 51 - FraxShare
Set-TextValue "D51" "6.96"
Set-TextValue "E51" "  +13.47%  "
